$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19").Insert(-4121, 0)

Write-Host "done"
